$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (the "date" column) switches from a numeric Date format to a
# plain Text format, and every existing value becomes a text string like
# "2023-03-30" instead of a date serial number. Apply the text format to
# the header + all data rows first so the subsequent string assignments
# are not re-interpreted as dates.
$dRange = $ws.Range("D1:D29")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "2023-03-30"
$ws.Range("D3").Value = "2023-03-31"
$ws.Range("D4").Value = "2023-04-01"
$ws.Range("D5").Value = "2023-04-02"
$ws.Range("D6").Value = "2023-04-03"
$ws.Range("D7").Value = "2023-04-04"
$ws.Range("D8").Value = "2023-04-05"
$ws.Range("D9").Value = "2023-04-06"
$ws.Range("D10").Value = "2023-04-07"
$ws.Range("D11").Value = "2023-04-08"
$ws.Range("D12").Value = "2023-04-09"
$ws.Range("D13").Value = "2023-04-10"
$ws.Range("D14").Value = "2023-04-02"
$ws.Range("D15").Value = "2023-04-03"
$ws.Range("D16").Value = "2023-04-04"
$ws.Range("D17").Value = "2023-04-02"
$ws.Range("D18").Value = "2023-04-03"
$ws.Range("D19").Value = "2023-04-04"
$ws.Range("D20").Value = "2023-04-03"
$ws.Range("D21").Value = "2023-04-04"
$ws.Range("D22").Value = "2023-04-05"
$ws.Range("D23").Value = "2023-04-06"
$ws.Range("D24").Value = "2023-04-07"
$ws.Range("D25").Value = "2023-04-08"
$ws.Range("D26").Value = "2023-04-09"
$ws.Range("D27").Value = "2023-04-10"
$ws.Range("D28").Value = "2023-04-04"
$ws.Range("D29").Value = "2023-04-06"

# Update the view: active cell/selection moves to D29 and the visible
# window scrolls up a couple of rows.
$ws.Activate()
$excel.Goto($ws.Range("A19"), $true) | Out-Null
$ws.Range("D29").Select() | Out-Null
